$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update runtime formulas for row 24 (Barts/Tahiti-column C, Tesla-column D) with newly collected data
$ws.Range("C24").Formula = "=(0.60382+0.61476+0.604918+0.73714)/4"
$ws.Range("D24").Formula = "=(0.818113+0.823809+0.819304+0.818177)/4"

# Update runtime formulas for row 25 with newly collected data
$ws.Range("C25").Formula = "=(3.77988+3.71921+3.76059+4.03667 )/4"
$ws.Range("D25").Formula = "=(2.22981+2.22904+2.22949+2.22917)/4"

# Add new formulas for row 27 (previously empty columns C and D)
$ws.Range("C27").Formula = "=(0.38893+0.428742+0.411657+0.411178)/4"
$ws.Range("D27").Formula = "=(0.556837+0.555766+0.559186+0.553763)/4"

# Add new formulas for row 28 (previously empty columns C and D)
$ws.Range("C28").Formula = "=(5.94859+5.16703+5.60233+5.71391)/4"
$ws.Range("D28").Formula = "=(3.48563+3.48459+3.48599+3.48248)/4"

# Update the active cell selection to reflect the author's final cursor position
[void]$ws.Range("G7").Select()
